$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123, shifting existing rows 123-187 down to 124-188.
$ws.Rows("123:123").Insert()

# Populate the newly inserted row 123 with the new price-record data.
$ws.Range("A123").Value = 3
$ws.Range("B123").Value = "Femacal de La Calera"
$ws.Range("C123").Value = "Coquimbo"
$ws.Range("D123").Value = 44529
$ws.Range("E123").Value = 5
$ws.Range("F123").Value = 100112001
$ws.Range("G123").Value = "Berenjena"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 85
$ws.Range("K123").Value = 7000
$ws.Range("L123").Value = 7500
$ws.Range("M123").Value = 7265
$ws.Range("N123").Value = "`$/caja 60 unidades"
$ws.Range("O123").Value = "Región de Arica y Parinacota"
$ws.Range("P123").Value = 121
$ws.Range("Q123").Value = 60
$ws.Range("R123").Value = "Hortaliza"
